$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (Price column updates, plus two E column text fixes)
$changes = @{
    "D2"  = "248.84"
    "D3"  = "22.53"
    "D4"  = "5.250"
    "D5"  = "0.05692"
    "D6"  = "3.411"
    "D7"  = "6.335"
    "D8"  = "0.8080"
    "D9"  = "0.8879"
    "D11" = "0.07446"
    "D12" = "0.03057"
    "D13" = "0.03102"
    "D14" = "0.09397"
    "D15" = "3.870"
    "D16" = "0.001585"
    "D17" = "0.04775"
    "D18" = "0.01828"
    "D19" = "0.0005810"
    "E19" = "18OneONE"
    "D20" = "0.006442"
    "D21" = "0.004983"
    "D22" = "0.0009963"
    "D24" = "3.688"
    "D25" = "2.167"
    "D27" = "0.1370"
    "D41" = "0.006840"
    "D42" = "0.1069"
    "D43" = "0.002808"
    "D44" = "0.007818"
    "D45" = "0.00005569"
    "D47" = "0.4990"
    "D48" = "0.2001"
    "E48" = "47BOLOBOLOWorstin24h"
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
}
